$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 20 - "Sample File": code block "block in all" -> "block in all #block
# stuff by default", explicit 24pt Courier size on the whole rule listing and
# drop the normAutofit fontScale override (back to plain <a:normAutofit/>).
# ---------------------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$sh20 = $s20.Shapes.Item(2)

$sh20.TextFrame.AutoSize = 2

$tr20 = $sh20.TextFrame.TextRange
$para20_1 = $tr20.Paragraphs(1)
$split20 = $para20_1.Characters(10, 3)
$split20.Text = "all #block stuff by default"

$tr20b = $sh20.TextFrame.TextRange
for ($i = 1; $i -le 4; $i++) {
    $tr20b.Paragraphs($i).Font.Size = 24
}

# ---------------------------------------------------------------------------
# Slide 21 - "Sample Rules 2": same text split, font size already 20pt so no
# explicit resize is required (and none happened in the source edit either).
# ---------------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$sh21 = $s21.Shapes.Item(2)

$tr21 = $sh21.TextFrame.TextRange
$para21_1 = $tr21.Paragraphs(1)
$split21 = $para21_1.Characters(10, 3)
$split21.Text = "all #block stuff by default"

# ---------------------------------------------------------------------------
# Slide 23 - "Macros cont'd": same text split on the 2nd paragraph (the first
# paragraph is the "Macros are then called..." text), font already 20pt.
# ---------------------------------------------------------------------------
$s23 = $p.Slides.Item(23)
$sh23 = $s23.Shapes.Item(2)

$tr23 = $sh23.TextFrame.TextRange
$para23_2 = $tr23.Paragraphs(2)
$split23 = $para23_2.Characters(10, 3)
$split23.Text = "all #block stuff by default"
